# fix: link to new variable id for freq demo data (#77)
#
# The freq demo table ("Tableau3") pointed rows 2-16 at the old/retired
# variable id "ser_pub_loc___variable_13". This relinks those rows to the
# new variable id "pop_region__population_totale" and appends 4 new rows
# of frequency data for a second new variable, "pop_region__type_region".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Relink the existing 15 rows (2-16) to the new variable id ------------
$ws.Range("A2:A16").Value = "pop_region__population_totale"

# --- Grow the table by 4 rows, then fill them in ---------------------------
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

$ws.Range("A17").Value = "pop_region__type_region"
$ws.Range("B17").Value = "urbaine"
$ws.Range("C17").Value = 9432

$ws.Range("A18").Value = "pop_region__type_region"
$ws.Range("B18").Value = "périurbaine"
$ws.Range("C18").Value = 43

$ws.Range("A19").Value = "pop_region__type_region"
$ws.Range("B19").Value = "rurale"
$ws.Range("C19").Value = 3434

$ws.Range("A20").Value = "pop_region__type_region"
$ws.Range("B20").Value = "montagne"
$ws.Range("C20").Value = 9481

# --- Widen columns A & B so the longer labels fit (mirrors author's autofit)
$ws.Columns.Item(1).ColumnWidth = 24.8
$ws.Columns.Item(2).ColumnWidth = 9

# --- Match the saved selection in the source file --------------------------
$ws.Range("A22").Select()
